$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: update heading text ---
$ws.Range("A1").Value = "posizione: (Università di Bologna, day 30, hour 12"

# --- Row 3 header columns (C:G) get new labels; drop column G ---
$ws.Range("C3").Value = "S_B"
$ws.Range("D3").Value = "S_D"
$ws.Range("E3").Value = "S"
$ws.Range("F3").Value = "% rispetto al max"
$ws.Range("G3").ClearContents()

# --- New data grid rows 4-15, columns A & B ---
$data = @(
    @(30, 0),
    @(30, 15),
    @(30, 30),
    @(15, 0),
    @(15, 15),
    @(15, 30),
    @(0, 0),
    @(0, 15),
    @(0, 30),
    @(90, 0),
    @(90, 15),
    @(90, 30)
)

$r = 4
foreach ($pair in $data) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r++
}

# --- Update selection to match the new active cell ---
$ws.Range("C4").Select()
